# Generate Report for Handback
# - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   on the Overview sheet (zh-cn/de-de columns) and on each language sheet's
#   Status column.
# - Each language sheet's "Latest Handback DateTime" is refreshed to the
#   handback-generation timestamp.
# - The stale "handback file is not the latest" error is cleared now that
#   the handback is in sync.
# - The two "...DateTime" report columns (Overview E/F, and each language
#   sheet's Status column C) are widened to fit the longer status text, and
#   the Error Detail column is narrowed now that it is typically empty.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

# ---- Overview sheet ------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Cells.Item(2, 5).Value = $statusNew   # E2 (zh-cn status)
$overview.Cells.Item(2, 6).Value = $statusNew   # F2 (de-de status)
$overview.Columns.Item(5).ColumnWidth = 29.17
$overview.Columns.Item(6).ColumnWidth = 29.17

# ---- zh-cn sheet -----------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Cells.Item(2, 3).Value = $statusNew                 # C2 Status
$zhcn.Cells.Item(2, 11).Value = "2016-09-01 22:54:17"      # K2 Latest Handback DateTime
$zhcn.Cells.Item(2, 16).Value = ""                         # P2 Error Detail cleared
$zhcn.Columns.Item(3).ColumnWidth = 29.17
$zhcn.Columns.Item(16).ColumnWidth = 12.8

# ---- de-de sheet -----------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Cells.Item(2, 3).Value = $statusNew                  # C2 Status
$dede.Cells.Item(2, 11).Value = "2016-09-01 22:54:24"      # K2 Latest Handback DateTime
$dede.Cells.Item(2, 16).Value = ""                         # P2 Error Detail cleared
$dede.Columns.Item(3).ColumnWidth = 29.17
$dede.Columns.Item(16).ColumnWidth = 12.8
